# GSC export update: append the 2025-11-10 row to the "Chart" data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Write the new date as literal text (not an auto-converted date serial) by
# building it through a text formula in a scratch cell and pasting the
# computed VALUE into A36 - this matches how the existing date-like labels
# in column A are stored (shared string, default/General style) instead of
# letting Excel reinterpret "2025-11-10" as a date value/format.
$ws.Range("Z1").Formula = "=""2025-11-10"""
$ws.Range("Z1").Copy()
$ws.Range("A36").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

$ws.Range("B36").Value = 0
$ws.Range("C36").Value = 66
